$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 22.75000000000012
$ws.Range("H2").Value = [double]"1.558207753859869e-16"
$ws.Range("I2").Value = 0.470267476287734
$ws.Range("K2").Value = 49.63182065922913
$ws.Range("L2").Value = "[46.7509407794014, 52.51270053905687]"
$ws.Range("O2").Value = 1.578658170272348
$ws.Range("P2").Value = "[1.515763422452732, 1.6415529180919641]"
$ws.Range("S2").Value = 53.11879401434737
$ws.Range("T2").Value = "[51.12561397402223, 55.111974054672515]"
$ws.Range("W2").Value = 17.03403403403412
$ws.Range("X2").Value = 16.80630630630639
$ws.Range("Y2").Value = 17.26176176176185

$ws.Range("E3").Value = 25.8300000000006
$ws.Range("H3").Value = [double]"1.558207753859869e-16"
$ws.Range("K3").Value = 45.70063980698307
$ws.Range("L3").Value = "[42.50055507014571, 48.90072454382043]"
$ws.Range("O3").Value = -0.08805264694746207
$ws.Range("P3").Value = "[-0.1635263443310011, -0.012578949563923025]"
$ws.Range("Q3").Value = 0.02237988394610957
$ws.Range("R3").Value = 0.02237988394610957
$ws.Range("S3").Value = 51.15535660246078
$ws.Range("T3").Value = "[49.060406113355334, 53.250307091566235]"
$ws.Range("W3").Value = 0.3619819819819945
$ws.Range("X3").Value = 0.05171171171171512
$ws.Range("Y3").Value = 0.6722522522522738
